$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column (D) keeps its original text formatting (e.g. "1.00",
# "0.0900", "13.55") instead of being auto-converted to a floating point number
# by the COM Value setter. Percent strings in column E and names/links in B/C are
# never auto-numeric so they do not need this treatment.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.939.24'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.290.61'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.93'
$ws.Range('E5').Value = '  -3.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.79'
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.616'
$ws.Range('E7').Value = '  -2.21%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.603'
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.97'
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0900'
$ws.Range('E11').Value = '  -1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.19'
$ws.Range('E12').Value = '  -4.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.107'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.980'
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.30'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.637.89'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.309.95'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.749.84'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.29'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.55'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000104'
$ws.Range('E21').Value = '  -1.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.39'
$ws.Range('E22').Value = '  -0.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '269.01'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.39'
$ws.Range('E24').Value = '  -5.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.19'
$ws.Range('E25').Value = '  -2.25%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.25'
$ws.Range('E27').Value = '  +17.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.76'
$ws.Range('E28').Value = '  -1.78%  '
$ws.Range('E29').Value = '  -1.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.37'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.72'
$ws.Range('E31').Value = '  -6.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.48'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0851'
$ws.Range('E33').Value = '  -3.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.130'
$ws.Range('E34').Value = '  -1.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.57'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('E36').Value = '  -2.80%  '
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('E38').Value = '  -2.78%  '
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.62'
$ws.Range('E40').Value = '  -3.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.00'
$ws.Range('E41').Value = '  +8.12%  '
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.87'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.227'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.01'
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('E46').Value = '  -3.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.733.98'
$ws.Range('E47').Value = '  +8.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '110.24'
$ws.Range('E48').Value = '  -3.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '76.94'
$ws.Range('E49').Value = '  -6.94%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.14'
$ws.Range('E50').Value = '  -2.99%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.62'
$ws.Range('E51').Value = '  -3.09%  '
